# Split the "Classifiers that have been used: ..." run into three runs,
# inserting a new "OPF, " run between "random forest, " and "KNN, Naïve Bayes".

$d = $word.ActiveDocument

# Locate the paragraph/run that needs to be split.
$target = $d.Content
$found = $target.Find.Execute(
    "Classifiers that have been used: CNN, SVM, random forest, KNN, Naïve Bayes",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the 'Classifiers that have been used...' paragraph text."
}

# Capture the paragraph's identity attributes so the rewritten paragraph
# keeps the same identity (paraId/rsid/etc.) instead of minting a new one.
$para = $target.Paragraphs.First
$paraXml = $para.Range.WordOpenXML
if ($paraXml -match '<w:p\s+([^>]*)>') {
    $paraAttrs = $matches[1]
} else {
    $paraAttrs = ""
}

# Clear the matched text (collapses the range to a single insertion point
# inside the now-empty paragraph) and rebuild it as three runs.
$target.Text = ""

$openXml = '<?xml version="1.0" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
    'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' +
    'xmlns:xml="http://www.w3.org/XML/1998/namespace">' +
    '<w:body>' +
    "<w:p $paraAttrs>" +
    '<w:r><w:t xml:space="preserve">Classifiers that have been used: CNN, SVM, random forest, </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">OPF, </w:t></w:r>' +
    '<w:r><w:t>KNN, Naïve Bayes</w:t></w:r>' +
    '</w:p>' +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($openXml) | Out-Null

# Sanity-check the result by re-finding the paragraph from scratch.
$check = $d.Content
$check.Find.Execute("Classifiers that have been used: CNN, SVM, random forest, OPF, KNN, Naïve Bayes",
                     $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Write-Output "Done: paragraph now reads -> $($check.Text)"
